# Swap the presentation's active theme palette from the "Integral" scheme
# to the "Office Theme" scheme (this mirrors the OOXML edit that exchanged
# the contents of ppt/theme/theme1.xml and ppt/theme/theme2.xml).
$p = $ppt.ActivePresentation

function HexToRGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

# Target palette ("Office Theme"), in the standard
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink order.
$officeTheme = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$cs = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $cs.Colors($i + 1).RGB = HexToRGBInt $officeTheme[$i]
}
